$wb = $excel.ActiveWorkbook

# --- Rename Task 1 sheet, add Task 2 sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Task 1 - Language Symbols"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Task 2 - Lexemes & RE"

# --- Populate Task 2 with the 26 "Letter = " lexeme class rows ---
$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 1
    $cell = $ws2.Cells.Item($row, 1)
    $cell.Value = $letters[$i] + " = "
}

Write-Host "done"
